$d = $word.ActiveDocument

$d.Content.Find.Execute("<id>p092v_a1</id>", $false, $false, $false, $false, $false,
                         $true, 1, $false, "<id>p092v_1</id>", 2)

$d.Content.Find.Execute("<id>p092v_a2</id>", $false, $false, $false, $false, $false,
                         $true, 1, $false, "<id>p092v_2</id>", 2)

$d.Content.Find.Execute("<id>p092v_a3</id>", $false, $false, $false, $false, $false,
                         $true, 1, $false, "<id>p092v_3</id>", 2)
